# Rename "Don't Know/Refused" -> "Doesn't Know/Prefers Not to Answer"
# across the data-quality issue names in column C of the EvaChecks sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C85").Value  = "Doesn't Know/Prefers Not to Answer Destination"
$ws.Range("C86").Value  = "Doesn't Know/Prefers Not to Answer Discharge Status"
$ws.Range("C87").Value  = "Doesn't Know/Prefers Not to Answer Gender"
$ws.Range("C88").Value  = "Doesn't Know/Prefers Not to Answer Length of Stay"
$ws.Range("C89").Value  = "Doesn't Know/Prefers Not to Answer Living Situation"
$ws.Range("C90").Value  = "Doesn't Know/Prefers Not to Answer Military Branch"
$ws.Range("C91").Value  = "Doesn't Know/Prefers Not to Answer Months or Times Homeless"
$ws.Range("C92").Value  = "Doesn't Know/Prefers Not to Answer Race/Ethnicity"
$ws.Range("C93").Value  = "Doesn't Know/Prefers Not to Answer Residence Prior"
$ws.Range("C94").Value  = "Doesn't Know/Prefers Not to Answer SSN"
$ws.Range("C95").Value  = "Doesn't Know/Prefers Not to Answer Veteran Status"
$ws.Range("C96").Value  = "Doesn't Know/Prefers Not to Answer War(s)"
$ws.Range("C97").Value  = "Doesn't Know/Prefers Not to Answer/Data Not Collected DOB"
$ws.Range("C103").Value = "Incomplete or Doesn't Know/Prefers Not to Answer Name"
